$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (Day 16)
$ws.Range("B13").Value = 30.36894366197182
$ws.Range("C13").Value = 26.61549295774649
$ws.Range("D13").Value = 26.61549295774649
$ws.Range("E13").Value = 29.13345070422535

# Row 14 (Day 17)
$ws.Range("B14").Value = 29.84014084507041
$ws.Range("C14").Value = 25.47957746478873
$ws.Range("D14").Value = 25.47957746478873
$ws.Range("E14").Value = 27.65450704225353

# Row 15 (Day 18)
$ws.Range("B15").Value = 27.40859154929578
$ws.Range("C15").Value = 23.47887323943662
$ws.Range("D15").Value = 23.47887323943662
$ws.Range("E15").Value = 25.25387323943661

# Row 16 (Day 19)
$ws.Range("B16").Value = 4967.880212765958
$ws.Range("C16").Value = 24.1372340425532
$ws.Range("D16").Value = 24.1372340425532
$ws.Range("E16").Value = 25.64478723404255
